$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was GLD data) -> now NEM data, labeled Newmont Corporation / NEM
$ws.Range("B2").Value = "Newmont Corporation"
$ws.Range("C2").Value = "NEM"
$ws.Range("D2").Value = 90.86
$ws.Range("E2").Value = 57.7
$ws.Range("F2").Value = 0.15
$ws.Range("G2").Value = 60
$ws.Range("H2").Value = 76
$ws.Range("I2").Value = 70
$ws.Range("J2").Value = 83
$ws.Range("K2").Value = 60.7
$ws.Range("N2").Value = 49.16024380385575

# Row 3 (was GC=F data) -> now GLD data, labeled StreetTRACKS Gold Shares / GLD
$ws.Range("B3").Value = "StreetTRACKS Gold Shares"
$ws.Range("C3").Value = "GLD"
$ws.Range("D3").Value = 387.73
$ws.Range("E3").Value = 71.2
$ws.Range("F3").Value = -0.04
$ws.Range("G3").Value = 40
$ws.Range("H3").Value = 76
$ws.Range("I3").Value = 83
$ws.Range("J3").Value = 93
$ws.Range("K3").Value = 59.9
$ws.Range("M3").Value = "⛔ 관망하십시오."
$ws.Range("N3").Value = 49.16024380385575

# Row 4 (was NEM data) -> now GC=F data, labeled Gold Feb 26 / GC=F
$ws.Range("B4").Value = "Gold Feb 26"
$ws.Range("C4").Value = "GC=F"
$ws.Range("D4").Value = 4239
$ws.Range("E4").Value = 71.5
$ws.Range("F4").Value = 0.49
$ws.Range("G4").Value = 50
$ws.Range("H4").Value = 70
$ws.Range("I4").Value = 73
$ws.Range("J4").Value = 80
$ws.Range("K4").Value = 58.9
$ws.Range("N4").Value = 49.16024380385575
